$d = $word.ActiveDocument

# Locate the target paragraph and rewrite its text.
$rng = $d.Content
$found = $rng.Find.Execute("Reify Kinds as SPOs.", $false, $true, $false, $false, $false, $true, 1, $false, "TBD: Relationship / Relation", 2)

# $rng now covers the replaced text; grab its paragraph so we can anchor
# the new paragraphs right after it.
$para = $rng.Paragraphs(1)

# Insert four new (initially empty) paragraphs right after it, in order.
$para.Range.InsertParagraphAfter()
$para.Range.InsertParagraphAfter()
$para.Range.InsertParagraphAfter()
$para.Range.InsertParagraphAfter()

# Fill in the text for the 2nd and 4th of the newly inserted paragraphs,
# leaving the 1st and 3rd (and the final, 5th) blank.
$para.Next().Next().Range.Text = "Reify Kinds as SPOs : Types Model"
$para.Next().Next().Next().Next().Range.Text = "Reify Statements as / Kinds / SPOs : Mappings Model"
